# Apply recomputed strike-count values (column G, header "K") to column G
# Commit message: "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals" -- the underlying K (strike) values were regenerated
# by an external data pipeline; here we write the resulting values directly
# into worksheet column G, matching the authoritative diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @{
    2 = 2
    3 = 0
    4 = 0
    5 = 1
    6 = 1
    7 = 1
    8 = 1
    9 = 0
    10 = 1
    11 = 1
    12 = 2
    13 = 2
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 2
    19 = 1
    20 = 0
    21 = 1
    22 = 1
    23 = 1
    24 = 2
    25 = 2
    26 = 2
    27 = 0
    28 = 1
    30 = 0
    31 = 2
    32 = 0
    33 = 0
    34 = 1
    35 = 1
    36 = 2
    37 = 1
    38 = 1
    39 = 3
    40 = 0
    41 = 0
    42 = 1
    43 = 0
    44 = 1
    45 = 0
    46 = 1
    47 = 1
    48 = 0
    49 = 1
    50 = 0
    51 = 0
    52 = 1
    53 = 1
    54 = 2
    55 = 3
    56 = 1
    57 = 0
    58 = 0
    59 = 2
    60 = 0
    61 = 0
    62 = 2
    63 = 1
    64 = 2
    65 = 0
    66 = 2
    67 = 0
    68 = 0
    69 = 1
    70 = 3
    71 = 1
    72 = 1
    73 = 1
    74 = 1
    75 = 2
    76 = 1
    77 = 1
    79 = 1
    80 = 1
    81 = 2
    82 = 1
    83 = 2
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $gValues[$row]
}
